# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (column G) previously held a placeholder/derived "Strike#" value.
# It is regenerated here to hold the actual strikeout total (K) for each outing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 2
    23 = 1
    24 = 3
    25 = 0
    26 = 3
    27 = 3
    28 = 1
    29 = 2
    30 = 2
    31 = 3
    32 = 1
    33 = 1
    34 = 1
    36 = 2
    38 = 2
    39 = 2
    40 = 0
    41 = 2
    42 = 0
    43 = 0
    44 = 1
    45 = 3
    46 = 3
    47 = 3
    48 = 1
    49 = 2
    50 = 1
    51 = 3
    52 = 1
    53 = 0
    54 = 2
    56 = 3
    57 = 5
    58 = 5
    59 = 2
    60 = 3
    61 = 1
    62 = 2
    63 = 1
    64 = 1
    65 = 2
    66 = 2
    67 = 0
    68 = 1
    69 = 3
    70 = 1
    72 = 2
    73 = 1
    74 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
